$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4
$ws.Range("B4").Value = "-"
$ws.Range("E4").Value = "MCT-3A-Manutenção mecânica"

# Row 6
$ws.Range("B6").Value = "-"
$ws.Range("E6").Value = "MCT-3A-Manutenção mecânica"

# Row 7
$ws.Range("B7").Value = "-"
$ws.Range("E7").Value = "MCT-3A-Manutenção mecânica"

# Row 8
$ws.Range("B8").Value = "-"
$ws.Range("E8").Value = "MCT-3A-Manutenção mecânica"

# Row 11
$ws.Range("B11").Value = "MEC-3A-Mec. Manut. Equip. Ind."
$ws.Range("E11").Value = "-"

# Row 12
$ws.Range("B12").Value = "MEC-3A-Mec. Manut. Equip. Ind."
$ws.Range("C12").Value = "-"
$ws.Range("E12").Value = "MEC-3A-Retificação"

# Row 14
$ws.Range("B14").Value = "MEC-3A-Retificação"
$ws.Range("C14").Value = "-"
$ws.Range("E14").Value = "MEC-3A-Mec. Manut. Equip. Ind."

# Row 15
$ws.Range("B15").Value = "MEC-3A-Retificação"
$ws.Range("C15").Value = "-"
$ws.Range("E15").Value = "MEC-3A-Mec. Manut. Equip. Ind."

# Row 16
$ws.Range("B16").Value = "MEC-3A-Retificação"
$ws.Range("C16").Value = "-"
